# Apply Freshservice ticket-data export + Meraki AP uptime refresh.

$wb = $excel.ActiveWorkbook

# --- 1. Update the "Meraki AP" sheet: new uptime readings + offline status ---
$ws = $wb.Worksheets.Item("Meraki AP")

$ws.Range("C2").Value  = 66.84
$ws.Range("D2").Value  = "offline"

$ws.Range("C3").Value  = 66.48999999999999
$ws.Range("D3").Value  = "offline"

$ws.Range("C11").Value = 67.44
$ws.Range("D11").Value = "offline"

$ws.Range("C12").Value = 67.43000000000001
$ws.Range("D12").Value = "offline"

$ws.Range("C14").Value = 67.43000000000001
$ws.Range("D14").Value = "offline"

$ws.Range("C15").Value = 67.33
$ws.Range("D15").Value = "offline"

$ws.Range("C16").Value = 67.43000000000001
$ws.Range("D16").Value = "offline"

$ws.Range("C17").Value = 67.44
$ws.Range("D17").Value = "offline"

$ws.Range("C18").Value = 67.44
$ws.Range("D18").Value = "offline"

$ws.Range("C19").Value = 67.44
$ws.Range("D19").Value = "offline"

$ws.Range("C20").Value = 67.44
$ws.Range("D20").Value = "offline"

$ws.Range("C21").Value = 67.33
$ws.Range("D21").Value = "offline"

$ws.Range("C22").Value = 67.33
$ws.Range("D22").Value = "offline"

$ws.Range("C23").Value = 67.43000000000001
$ws.Range("D23").Value = "offline"

$ws.Range("C24").Value = 67.44
$ws.Range("D24").Value = "offline"

$ws.Range("C25").Value = 67.43000000000001
$ws.Range("D25").Value = "offline"

$ws.Range("C26").Value = 67.44
$ws.Range("D26").Value = "offline"

# --- 2. Add a new "Freshservice" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$fs = $wb.Worksheets.Add($null, $lastSheet)
$fs.Name = "Freshservice"

$fs.Range("A1").Value = "Total Tickets (Last 7 Days)"
$fs.Range("B1").Value = "Unresolved Tickets (Last 7 Days)"
$fs.Range("C1").Value = "Resolved Tickets (Last 7 Days)"

$fs.Range("A2").Value = 27
$fs.Range("B2").Value = 7
$fs.Range("C2").Value = 20
